$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M40").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0

$ws.Range("M76").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0

$ws.Range("M79").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0

$ws.Range("H112").Value = 1496.3334
$ws.Range("I112").Value = 489.5
$ws.Range("J112").Value = 1999.75
$ws.Range("K112").Value = 1468.5
$ws.Range("L112").Value = 5999.25
$ws.Range("M112").Value = -360.5
$ws.Range("N112").Value = -8215.25

$ws.Range("H125").Value = 1508
$ws.Range("I125").Value = 1289.125
$ws.Range("J125").Value = 1799.8334
$ws.Range("K125").Value = 11602.125
$ws.Range("L125").Value = 16198.5006
$ws.Range("M125").Value = -9142.125
$ws.Range("N125").Value = -21118.5006

$ws.Range("H137").Value = 2099.2
$ws.Range("I137").Value = 1606.8334
$ws.Range("J137").Value = 2837.75
$ws.Range("K137").Value = 4820.5002
$ws.Range("L137").Value = 8513.25
$ws.Range("M137").Value = -2270.5002
$ws.Range("N137").Value = -13613.25

$ws.Range("H138").Value = 7469.75
$ws.Range("J138").Value = 7216.6924
$ws.Range("L138").Value = 21650.0772
$ws.Range("N138").Value = -31930.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4948
$ws.Range("I2").Value = 4948
$ws.Range("K2").Value = 4948
$ws.Range("M2").Value = -4835

$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -788

$ws.Range("H74").Value = 4012
$ws.Range("I74").Value = 4012
$ws.Range("K74").Value = 4012
$ws.Range("M74").Value = -3138

$ws.Range("H77").Value = 4012
$ws.Range("I77").Value = 4012
$ws.Range("K77").Value = 20060
$ws.Range("M77").Value = -15692

$ws.Range("H101").Value = 34800.668
$ws.Range("J101").Value = 34800.668
$ws.Range("L101").Value = 34800.668
$ws.Range("N101").Value = -41290.668

$ws.Range("H116").Value = 4948
$ws.Range("I116").Value = 4948
$ws.Range("K116").Value = 4948
$ws.Range("M116").Value = -2654

$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820

$ws.Range("H132").Value = 4166.3335
$ws.Range("I132").Value = 3249.5
$ws.Range("K132").Value = 9748.5
$ws.Range("M132").Value = -7218.5

$ws.Range("H135").Value = 38974.168
$ws.Range("J135").Value = 38974.168
$ws.Range("L135").Value = 38974.168
$ws.Range("N135").Value = -49114.168

$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4948
$ws.Range("I3").Value = 4948
$ws.Range("K3").Value = 4948
$ws.Range("M3").Value = -4834

$ws.Range("H22").Value = 495.7143
$ws.Range("I22").Value = 580.2727
$ws.Range("J22").Value = 185.66667
$ws.Range("K22").Value = 580.2727
$ws.Range("L22").Value = 185.66667
$ws.Range("M22").Value = -407.2727
$ws.Range("N22").Value = -531.6666700000001

$ws.Range("H134").Value = 3159.3333
$ws.Range("I134").Value = 3191.8
$ws.Range("K134").Value = 9575.400000000001
$ws.Range("M134").Value = -7040.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 246.07692
$ws.Range("I7").Value = 174.85715
$ws.Range("K7").Value = 174.85715
$ws.Range("M7").Value = -61.85714999999999

$ws.Range("H22").Value = 792.1429
$ws.Range("I22").Value = 807.5
$ws.Range("K22").Value = 807.5
$ws.Range("M22").Value = -457.5

$ws.Range("H41").Value = 30141.666
$ws.Range("I41").Value = 5359
$ws.Range("J41").Value = 33239.5
$ws.Range("K41").Value = 5359
$ws.Range("L41").Value = 33239.5
$ws.Range("M41").Value = -4931
$ws.Range("N41").Value = -34095.5

$ws.Range("H107").Value = 692.8333
$ws.Range("I107").Value = 351.875
$ws.Range("K107").Value = 351.875
$ws.Range("M107").Value = 1568.125

$ws.Range("N132").ClearContents()
$ws.Range("H132").Value = 200
$ws.Range("I132").Value = 200
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 600
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 1930

$ws.Range("N139").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H141").Value = 972457.5
$ws.Range("J141").Value = 972457.5
$ws.Range("L141").Value = 972457.5
$ws.Range("N141").Value = -982817.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9499.5
$ws.Range("I56").Value = 9499.5
$ws.Range("K56").Value = 9499.5
$ws.Range("M56").Value = -8969.5

$ws.Range("H107").Value = 771.8571
$ws.Range("J107").Value = 1916.5
$ws.Range("L107").Value = 5749.5
$ws.Range("N107").Value = -9589.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N132").ClearContents()
$ws.Range("H132").Value = 5001
$ws.Range("I132").Value = 5001
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15003
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12473

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 907.1818
$ws.Range("I16").Value = 885
$ws.Range("J16").Value = 966.3333
$ws.Range("K16").Value = 885
$ws.Range("L16").Value = 966.3333
$ws.Range("M16").Value = -715
$ws.Range("N16").Value = -1306.3333

$ws.Range("M132").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M28").ClearContents()
$ws.Range("H28").Value = 20000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 20000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 20000
$ws.Range("N28").Value = -20696

$ws.Range("H81").Value = 6250.5
$ws.Range("I81").Value = 6250.5
$ws.Range("K81").Value = 12501
$ws.Range("M81").Value = -11440

$ws.Range("H84").Value = 6250.5
$ws.Range("I84").Value = 6250.5
$ws.Range("K84").Value = 62505
$ws.Range("M84").Value = -57201
